$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H93").Value = 175397.75
$ws.Range("J93").Value = 175397.75
$ws.Range("L93").Value = 175397.75
$ws.Range("N93").Value = -180389.75
$ws.Range("H116").Value = 3849.1
$ws.Range("I116").Value = 3563.125
$ws.Range("K116").Value = 3563.125
$ws.Range("M116").Value = -121.125
$ws.Range("H125").Value = 2730
$ws.Range("I125").Value = 1987.5
$ws.Range("J125").Value = 3000
$ws.Range("K125").Value = 17887.5
$ws.Range("L125").Value = 27000
$ws.Range("M125").Value = -15427.5
$ws.Range("N125").Value = -31920

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1247984.9
$ws.Range("I32").Value = 1392666.2
$ws.Range("K32").Value = 1392666.2
$ws.Range("M32").Value = -1392379.2
$ws.Range("H102").Value = 3086.6667
$ws.Range("I102").Value = 1560
$ws.Range("J102").Value = 4995
$ws.Range("K102").Value = 1560
$ws.Range("L102").Value = 4995
$ws.Range("M102").Value = 62
$ws.Range("N102").Value = -8239
$ws.Range("H132").Value = 26199.455
$ws.Range("I132").Value = 40182.742
$ws.Range("J132").Value = 3990.7058
$ws.Range("K132").Value = 120548.226
$ws.Range("L132").Value = 11972.1174
$ws.Range("M132").Value = -118018.226
$ws.Range("N132").Value = -17032.1174

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1748.1177
$ws.Range("I94").Value = 661.75
$ws.Range("J94").Value = 2713.7778
$ws.Range("K94").Value = 661.75
$ws.Range("L94").Value = 2713.7778
$ws.Range("M94").Value = -210.75
$ws.Range("N94").Value = -3615.7778
$ws.Range("H103").Value = 39000
$ws.Range("J103").Value = 39000
$ws.Range("L103").Value = 39000
$ws.Range("N103").Value = -41344
$ws.Range("H105").Value = 1776.0769
$ws.Range("I105").Value = 1497.8
$ws.Range("J105").Value = 1950
$ws.Range("K105").Value = 1497.8
$ws.Range("L105").Value = 1950
$ws.Range("M105").Value = 249.2
$ws.Range("N105").Value = -5444
$ws.Range("H107").Value = 1426.2106
$ws.Range("I107").Value = 1061.3846
$ws.Range("J107").Value = 2216.6667
$ws.Range("K107").Value = 1061.3846
$ws.Range("L107").Value = 2216.6667
$ws.Range("M107").Value = 858.6153999999999
$ws.Range("N107").Value = -6056.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 750
$ws.Range("I16").Value = 750
$ws.Range("K16").Value = 750
$ws.Range("M16").Value = -463
$ws.Range("H94").Value = 5619.2354
$ws.Range("I94").Value = 833.3333
$ws.Range("J94").Value = 6644.7856
$ws.Range("K94").Value = 833.3333
$ws.Range("L94").Value = 6644.7856
$ws.Range("M94").Value = -382.3333
$ws.Range("N94").Value = -7546.7856
$ws.Range("H105").Value = 910.43475
$ws.Range("I105").Value = 735.625
$ws.Range("J105").Value = 1310
$ws.Range("K105").Value = 735.625
$ws.Range("L105").Value = 1310
$ws.Range("M105").Value = 1011.375
$ws.Range("N105").Value = -4804
$ws.Range("H113").Value = 750
$ws.Range("I113").Value = 750
$ws.Range("K113").Value = 750
$ws.Range("M113").Value = 1420
$ws.Range("H122").Value = 1522.238
$ws.Range("I122").Value = 935.38464
$ws.Range("J122").Value = 2475.875
$ws.Range("K122").Value = 2806.15392
$ws.Range("L122").Value = 7427.625
$ws.Range("M122").Value = -356.1539199999997
$ws.Range("N122").Value = -12327.625
$ws.Range("H132").Value = 2032.7838
$ws.Range("I132").Value = 1132.2307
$ws.Range("K132").Value = 3396.6921
$ws.Range("M132").Value = -866.6921000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 4755.579
$ws.Range("I116").Value = 664.5
$ws.Range("J116").Value = 5846.533
$ws.Range("K116").Value = 1993.5
$ws.Range("L116").Value = 17539.599
$ws.Range("M116").Value = 1448.5
$ws.Range("N116").Value = -24423.599
$ws.Range("H122").Value = 739.3125
$ws.Range("I122").Value = 305
$ws.Range("J122").Value = 1297.7142
$ws.Range("K122").Value = 2745
$ws.Range("L122").Value = 11679.4278
$ws.Range("M122").Value = -295
$ws.Range("N122").Value = -16579.4278
$ws.Range("H129").Value = 2571.68
$ws.Range("I129").Value = 4270
$ws.Range("J129").Value = 1772.4706
$ws.Range("K129").Value = 12810
$ws.Range("L129").Value = 5317.4118
$ws.Range("M129").Value = -7810
$ws.Range("N129").Value = -15317.4118
$ws.Range("H131").Value = 2607.5615
$ws.Range("I131").Value = 3383.8462
$ws.Range("J131").Value = 2378.2046
$ws.Range("K131").Value = 10151.5386
$ws.Range("L131").Value = 7134.6138
$ws.Range("M131").Value = -5111.5386
$ws.Range("N131").Value = -17214.6138

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3044.111
$ws.Range("I102").Value = 1470.96
$ws.Range("K102").Value = 1470.96
$ws.Range("M102").Value = 151.04
$ws.Range("H107").Value = 5204.35
$ws.Range("J107").Value = 346.7143
$ws.Range("L107").Value = 346.7143
$ws.Range("N107").Value = -4186.7143
$ws.Range("H132").Value = 4007.195
$ws.Range("I132").Value = 4051.8076
$ws.Range("J132").Value = 3929.8667
$ws.Range("K132").Value = 12155.4228
$ws.Range("L132").Value = 11789.6001
$ws.Range("M132").Value = -9625.4228
$ws.Range("N132").Value = -16849.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2747.682
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 2573
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 2573
$ws.Range("M7").Value = -2888
$ws.Range("N7").Value = -2797
$ws.Range("H16").Value = 1006.52
$ws.Range("I16").Value = 858.5333000000001
$ws.Range("J16").Value = 2338.4
$ws.Range("K16").Value = 858.5333000000001
$ws.Range("L16").Value = 2338.4
$ws.Range("M16").Value = -688.5333000000001
$ws.Range("N16").Value = -2678.4
$ws.Range("H55").Value = 167.91304
$ws.Range("I55").Value = 139.41176
$ws.Range("J55").Value = 248.66667
$ws.Range("K55").Value = 139.41176
$ws.Range("L55").Value = 248.66667
$ws.Range("M55").Value = 33.58824000000001
$ws.Range("N55").Value = -594.6666700000001
$ws.Range("H61").Value = 1764.9
$ws.Range("I61").Value = 851
$ws.Range("K61").Value = 851
$ws.Range("M61").Value = -649
$ws.Range("H93").Value = 946.2069
$ws.Range("I93").Value = 910.8333
$ws.Range("J93").Value = 1116
$ws.Range("K93").Value = 910.8333
$ws.Range("L93").Value = 1116
$ws.Range("M93").Value = 337.1667
$ws.Range("N93").Value = -3612
$ws.Range("H113").Value = 1764.9
$ws.Range("I113").Value = 851
$ws.Range("K113").Value = 851
$ws.Range("M113").Value = 1319
$ws.Range("H122").Value = 2534.5264
$ws.Range("I122").Value = 2531.0667
$ws.Range("K122").Value = 7593.2001
$ws.Range("M122").Value = -5143.2001
$ws.Range("H126").Value = 2747.682
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 2573
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 7719
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -12659

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 28857.5
$ws.Range("J92").Value = 28857.5
$ws.Range("L92").Value = 28857.5
$ws.Range("N92").Value = -33849.5
